$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.6537219968967349
$ws.Range("C2").Value2 = 0.07532276124402415
$ws.Range("D2").Value2 = 0.1200858742501794
$ws.Range("E2").Value2 = 0.1232652603590836
$ws.Range("F2").Value2 = 1.86117351845887
$ws.Range("I2").Value2 = 1.203488368999011
$ws.Range("J2").Value2 = 0.1600578772562233
$ws.Range("K2").Value2 = 0.406884506053018
$ws.Range("L2").Value2 = 0.2914885544889643
$ws.Range("O2").Value2 = 4.917389103029834

$ws.Range("B3").Value2 = 0.6139160630081903
$ws.Range("C3").Value2 = 0.07312015266788308
$ws.Range("D3").Value2 = 0.1180219903983684
$ws.Range("E3").Value2 = 0.1233192762499478
$ws.Range("F3").Value2 = 1.870935651369912
$ws.Range("I3").Value2 = 1.213442784878769
$ws.Range("J3").Value2 = 0.1610591200375726
$ws.Range("K3").Value2 = 0.3695605442178191
$ws.Range("L3").Value2 = 0.2870815846421166
$ws.Range("O3").Value2 = 4.951555263342044

$ws.Range("B4").Value2 = 0.5896291302998691
$ws.Range("C4").Value2 = 0.07175443171388451
$ws.Range("D4").Value2 = 0.1167940951050497
$ws.Range("E4").Value2 = 0.1233928815279963
$ws.Range("F4").Value2 = 1.877749120894762
$ws.Range("I4").Value2 = 1.220060256618243
$ws.Range("J4").Value2 = 0.1617227044180822
$ws.Range("K4").Value2 = 0.3466825284026811
$ws.Range("L4").Value2 = 0.2844843249489202
$ws.Range("O4").Value2 = 4.974739993357659

$ws.Range("B5").Value2 = 0.5797715301365258
$ws.Range("C5").Value2 = 0.07119456716172579
$ws.Range("D5").Value2 = 0.1163036776228097
$ws.Range("E5").Value2 = 0.1234330743918068
$ws.Range("F5").Value2 = 1.880731935628276
$ws.Range("I5").Value2 = 1.222884101321661
$ws.Range("J5").Value2 = 0.1620054149923913
$ws.Range("K5").Value2 = 0.3373700044115253
$ws.Range("L5").Value2 = 0.2834533730322377
$ws.Range("O5").Value2 = 4.984743170800769

$ws.Range("B6").Value2 = 0.5781370931603362
$ws.Range("C6").Value2 = 0.07110140216116889
$ws.Range("D6").Value2 = 0.1162228477250409
$ws.Range("E6").Value2 = 0.1234403651897757
$ws.Range("F6").Value2 = 1.881239694021922
$ws.Range("I6").Value2 = 1.223360681745188
$ws.Range("J6").Value2 = 0.1620531019889917
$ws.Range("K6").Value2 = 0.3358243169551827
$ws.Range("L6").Value2 = 0.2832838462805327
$ws.Range("O6").Value2 = 4.986437734573187

$ws.Range("B7").Value2 = 0.5894960261979634
$ws.Range("C7").Value2 = 0.07174689460251926
$ws.Range("D7").Value2 = 0.1167874407455969
$ws.Range("E7").Value2 = 0.1233933822521269
$ws.Range("F7").Value2 = 1.877788512723484
$ws.Range("I7").Value2 = 1.220097824959225
$ws.Range("J7").Value2 = 0.1617264673464067
$ws.Range("K7").Value2 = 0.3465568932800807
$ws.Range("L7").Value2 = 0.2844703098380705
$ws.Range("O7").Value2 = 4.974872651211641

$ws.Range("B8").Value2 = 0.639965387220883
$ws.Range("C8").Value2 = 0.07456607862346942
$ws.Range("D8").Value2 = 0.1193661224789935
$ws.Range("E8").Value2 = 0.1232755092399245
$ws.Range("F8").Value2 = 1.864369540681054
$ws.Range("I8").Value2 = 1.206815799451473
$ws.Range("J8").Value2 = 0.1603929851968005
$ws.Range("K8").Value2 = 0.3940074679550492
$ws.Range("L8").Value2 = 0.2899465623740412
$ws.Range("O8").Value2 = 4.928711831321124

$ws.Range("B9").Value2 = 0.7401291332955964
$ws.Range("C9").Value2 = 0.07998802989307308
$ws.Range("D9").Value2 = 0.1247323883476952
$ws.Range("E9").Value2 = 0.1233640259401234
$ws.Range("F9").Value2 = 1.844548893484294
$ws.Range("I9").Value2 = 1.184776582286858
$ws.Range("J9").Value2 = 0.1581645298034253
$ws.Range("K9").Value2 = 0.4873445318895619
$ws.Range("L9").Value2 = 0.301542335377377
$ws.Range("O9").Value2 = 4.855684377287588

$ws.Range("B10").Value2 = 0.814415748821375
$ws.Range("C10").Value2 = 0.08390585438466758
$ws.Range("D10").Value2 = 0.1288605913726428
$ws.Range("E10").Value2 = 0.123622445827106
$ws.Range("F10").Value2 = 1.833935473120349
$ws.Range("I10").Value2 = 1.171022323251808
$ws.Range("J10").Value2 = 0.1567617587040004
$ws.Range("K10").Value2 = 0.5560697233172505
$ws.Range("L10").Value2 = 0.3105783719935573
$ws.Range("O10").Value2 = 4.812678099511857

$ws.Range("B11").Value2 = 0.8483557150259173
$ws.Range("C11").Value2 = 0.08567376906763968
$ws.Range("D11").Value2 = 0.1307783039193993
$ws.Range("E11").Value2 = 0.1237816927376478
$ws.Range("F11").Value2 = 1.829962555899712
$ws.Range("I11").Value2 = 1.165293556979393
$ws.Range("J11").Value2 = 0.1561742746763812
$ws.Range("K11").Value2 = 0.5873625494513419
$ws.Range("L11").Value2 = 0.3148001339626774
$ws.Range("O11").Value2 = 4.795421749046483

$ws.Range("B12").Value2 = 0.8612282962815527
$ws.Range("C12").Value2 = 0.08634115256617747
$ws.Range("D12").Value2 = 0.1315101488158916
$ws.Range("E12").Value2 = 0.1238479627023281
$ws.Range("F12").Value2 = 1.828580907772917
$ws.Range("I12").Value2 = 1.16320009323583
$ws.Range("J12").Value2 = 0.1559590734988916
$ws.Range("K12").Value2 = 0.5992160067431769
$ws.Range("L12").Value2 = 0.3164146788525528
$ws.Range("O12").Value2 = 4.789218726655463

$ws.Range("B13").Value2 = 0.8584550680374718
$ws.Range("C13").Value2 = 0.08619751286809674
$ws.Range("D13").Value2 = 0.1313522828127418
$ws.Range("E13").Value2 = 0.1238334253649533
$ws.Range("F13").Value2 = 1.828873011059784
$ws.Range("I13").Value2 = 1.163647583367649
$ws.Range("J13").Value2 = 0.1560050980104783
$ws.Range("K13").Value2 = 0.596663005905782
$ws.Range("L13").Value2 = 0.3160662547606989
$ws.Range("O13").Value2 = 4.790539914794579

$ws.Range("B14").Value2 = 0.8494143488860573
$ws.Range("C14").Value2 = 0.08572871710232732
$ws.Range("D14").Value2 = 0.1308384004326371
$ws.Range("E14").Value2 = 0.1237870253902393
$ws.Range("F14").Value2 = 1.829846426536875
$ws.Range("I14").Value2 = 1.165119805660332
$ws.Range("J14").Value2 = 0.156156424411467
$ws.Range("K14").Value2 = 0.5883376734130934
$ws.Range("L14").Value2 = 0.3149326465727142
$ws.Range("O14").Value2 = 4.794904778511352

$ws.Range("B15").Value2 = 0.8438792602580918
$ws.Range("C15").Value2 = 0.08544129349513696
$ws.Range("D15").Value2 = 0.1305243661863216
$ws.Range("E15").Value2 = 0.1237593802530412
$ws.Range("F15").Value2 = 1.830458660144743
$ws.Range("I15").Value2 = 1.166031466806814
$ws.Range("J15").Value2 = 0.1562500620048226
$ws.Range("K15").Value2 = 0.5832386080516017
$ws.Range("L15").Value2 = 0.3142403395370934
$ws.Range("O15").Value2 = 4.797621558899152

$ws.Range("B16").Value2 = 0.8122005714181739
$ws.Range("C16").Value2 = 0.08379002640641886
$ws.Range("D16").Value2 = 0.1287360586440514
$ws.Range("E16").Value2 = 0.1236128751546843
$ws.Range("F16").Value2 = 1.834212306966947
$ws.Range("I16").Value2 = 1.171407334638968
$ws.Range("J16").Value2 = 0.1568011698579959
$ws.Range("K16").Value2 = 0.5540251965164771
$ws.Range("L16").Value2 = 0.3103046961762459
$ws.Range("O16").Value2 = 4.813852254889099

$ws.Range("B17").Value2 = 0.792803681486788
$ws.Range("C17").Value2 = 0.08277333919741636
$ws.Range("D17").Value2 = 0.1276491263746493
$ws.Range("E17").Value2 = 0.1235336567039873
$ws.Range("F17").Value2 = 1.836733963027569
$ws.Range("I17").Value2 = 1.174840485164538
$ws.Range("J17").Value2 = 0.1571522161507168
$ws.Range("K17").Value2 = 0.5361107558922242
$ws.Range("L17").Value2 = 0.3079186948495476
$ws.Range("O17").Value2 = 4.824400089573288

$ws.Range("B18").Value2 = 0.7816609527292826
$ws.Range("C18").Value2 = 0.08218722107596932
$ws.Range("D18").Value2 = 0.1270276994397648
$ws.Range("E18").Value2 = 0.1234920181411177
$ws.Range("F18").Value2 = 1.838264856423407
$ws.Range("I18").Value2 = 1.176864854010525
$ws.Range("J18").Value2 = 0.1573588964717167
$ws.Range("K18").Value2 = 0.5258096506865968
$ws.Range("L18").Value2 = 0.3065568071608453
$ws.Range("O18").Value2 = 4.830684123705197

$ws.Range("B19").Value2 = 0.777890620190135
$ws.Range("C19").Value2 = 0.08198854112987419
$ws.Range("D19").Value2 = 0.1268179407858838
$ws.Range("E19").Value2 = 0.1234785951648263
$ws.Range("F19").Value2 = 1.838797022194669
$ws.Range("I19").Value2 = 1.177558809935878
$ws.Range("J19").Value2 = 0.1574296942425129
$ws.Range("K19").Value2 = 0.5223223764581917
$ws.Range("L19").Value2 = 0.3060974990769552
$ws.Range("O19").Value2 = 4.832849100508383

$ws.Range("B20").Value2 = 0.7948670862027427
$ws.Range("C20").Value2 = 0.08288170684777185
$ws.Range("D20").Value2 = 0.1277644448018407
$ws.Range("E20").Value2 = 0.1235416835443068
$ws.Range("F20").Value2 = 1.836457197813054
$ws.Range("I20").Value2 = 1.174469875817241
$ws.Range("J20").Value2 = 0.1571143533441699
$ws.Range("K20").Value2 = 0.5380174930984651
$ws.Range("L20").Value2 = 0.3081716054210801
$ws.Range("O20").Value2 = 4.823254776923278

$ws.Range("B21").Value2 = 0.8520692842037363
$ws.Range("C21").Value2 = 0.08586647065540376
$ws.Range("D21").Value2 = 0.1309891873990807
$ws.Range("E21").Value2 = 0.1238004924816316
$ws.Range("F21").Value2 = 1.82955717911608
$ws.Range("I21").Value2 = 1.164685318868294
$ws.Range("J21").Value2 = 0.1561117791036502
$ws.Range("K21").Value2 = 0.5907829346509459
$ws.Range("L21").Value2 = 0.3152651855729971
$ws.Range("O21").Value2 = 4.793613714927972

$ws.Range("B22").Value2 = 0.8895718917054296
$ws.Range("C22").Value2 = 0.08780500647948486
$ws.Range("D22").Value2 = 0.1331296358850977
$ws.Range("E22").Value2 = 0.1240044039316111
$ws.Range("F22").Value2 = 1.825763368683042
$ws.Range("I22").Value2 = 1.158732902097
$ws.Range("J22").Value2 = 0.1554988848211778
$ws.Range("K22").Value2 = 0.6252885375755
$ws.Range("L22").Value2 = 0.3199936026251038
$ws.Range("O22").Value2 = 4.776174184310094

$ws.Range("B23").Value2 = 0.8695455516364916
$ws.Range("C23").Value2 = 0.08677149701617282
$ws.Range("D23").Value2 = 0.1319842517930141
$ws.Range("E23").Value2 = 0.1238924009436282
$ws.Range("F23").Value2 = 1.82772275990645
$ws.Range("I23").Value2 = 1.161869358618475
$ws.Range("J23").Value2 = 0.1558221287133925
$ws.Range("K23").Value2 = 0.6068706123786001
$ws.Range("L23").Value2 = 0.3174615529569991
$ws.Range("O23").Value2 = 4.785305228566614

$ws.Range("B24").Value2 = 0.7939341934524862
$ws.Range("C24").Value2 = 0.08283271885068189
$ws.Range("D24").Value2 = 0.1277122985430736
$ws.Range("E24").Value2 = 0.123538042443613
$ws.Range("F24").Value2 = 1.83658207054696
$ws.Range("I24").Value2 = 1.174637270668079
$ws.Range("J24").Value2 = 0.1571314559823271
$ws.Range("K24").Value2 = 0.5371554629337822
$ws.Range("L24").Value2 = 0.3080572338407137
$ws.Range("O24").Value2 = 4.82377188764309

$ws.Range("B25").Value2 = 0.7129076574818214
$ws.Range("C25").Value2 = 0.0785327271242906
$ws.Range("D25").Value2 = 0.1232478561651078
$ws.Range("E25").Value2 = 0.1233060116639102
$ws.Range("F25").Value2 = 1.849216724904643
$ws.Range("I25").Value2 = 1.190310306174013
$ws.Range("J25").Value2 = 0.1587261273452043
$ws.Range("K25").Value2 = 0.4620662161180462
$ws.Range("L25").Value2 = 0.2983142331406583
$ws.Range("O25").Value2 = 4.873569239585322

